# Update countries & provincias Spain
# - Reorder Montserrat/Seychelles (Montserrat now listed before Seychelles)
#   and carry each country's statistics row with it.
# - Update the "Datos actualizados" timestamp footer.
# - Refresh the numeric statistics (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes hoy, Muertes) for several countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 21:03"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1252204
$ws.Range("C4").Value = 14571
$ws.Range("D4").Value = 204657
$ws.Range("E4").Value = 974050
$ws.Range("F4").Value = 15990
$ws.Range("G4").Value = 1226
$ws.Range("H4").Value = 73497

# --- Francia (row 8) ---
$ws.Range("B8").Value = 174191
$ws.Range("C8").Value = 3640
$ws.Range("E8").Value = 94410

# --- Alemania (row 9) ---
$ws.Range("B9").Value = 167575
$ws.Range("C9").Value = 568
$ws.Range("E9").Value = 22985
$ws.Range("G9").Value = 197
$ws.Range("H9").Value = 7190

# --- India (row 16) ---
$ws.Range("B16").Value = 52987
$ws.Range("C16").Value = 3587
$ws.Range("D16").Value = 15331
$ws.Range("E16").Value = 35871

# --- Kazajistan (row 60) ---
$ws.Range("B60").Value = 4422
$ws.Range("C60").Value = 217
$ws.Range("E60").Value = 2984

# --- Estado de Palestina (row 128) ---
$ws.Range("D128").Value = 174
$ws.Range("E128").Value = 195

# --- Isla de Man (row 131) ---
$ws.Range("B131").Value = 327
$ws.Range("C131").Value = 1
$ws.Range("E131").Value = 33

# --- Monaco (row 161) ---
$ws.Range("D161").Value = 82
$ws.Range("E161").Value = 9

# --- Swap Montserrat / Seychelles (rows 205-206), carrying each country's
#     own statistics with its new label ---
$ws.Range("A205").Value = "Montserrat"
$ws.Range("D205").Value = 7
$ws.Range("F205").Value = 1
$ws.Range("H205").Value = 1

$ws.Range("A206").Value = "Seychelles"
$ws.Range("D206").Value = 8
$ws.Range("F206").Value = 0
$ws.Range("H206").Value = 0
